$d = $word.ActiveDocument

# Target the first paragraph (the hidden ID marker paragraph).
$p1 = $d.Paragraphs.Item(1)

# Replace the two runs ("**ID__AFFARS_pgi_5309_topic_3__ID**" + " ") with a
# single run containing the updated marker text (no trailing space).
$p1.Range.Find.Execute("**ID__AFFARS_pgi_5309_topic_3__ID** ", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "**ID__AFFARS_5309_104_6__ID**", 2)

# Add a 4-sided paragraph border (5pt space from text on each side).
$p1.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

Write-Output ("Paragraph 1 now reads: " + $d.Paragraphs.Item(1).Range.Text)
